$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("soknad")

# Insert three new rows above the current row 2, pushing the existing
# applications (old rows 2-7) down to rows 5-10. This preserves the
# per-application data (foresatt_1/2, barn_1, barnehager_prioritert,
# tidspunkt_oppstart, brutto_inntekt, etc.) together as it moves.
$ws.Rows("2:4").Insert()

# The newly inserted rows inherited bold/border formatting from the row
# above; strip that back to the plain (unstyled) look used by blank data
# rows elsewhere in this sheet.
$ws.Range("B2:M4").ClearFormats()

# Column A (the index column) keeps the bold/centered/bordered style used
# throughout the sheet; carry that style onto the freshly inserted cells
# before filling them in.
$ws.Range("A5").Copy($ws.Range("A2:A4"))

# Renumber the leading index column (A) sequentially for every data row.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8

# sok_id (column B) becomes a descending run from 9 down to 1 across the
# now-10 data rows; rows 5-10 already carry the correct shifted values,
# so only the three freshly inserted rows need to be filled in.
$ws.Range("B2").Value = 9
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = 7
